$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the XML-mapped table from A1:I4 to A1:K4 (adds 2 new table columns).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K4"))

# Header text for the two new columns. Write K1 (ActorID) before J1 (CanClone)
# so the shared-string table gets "ActorID" then "CanClone" appended in that
# order (matches the target shared strings: ...,newscene,ActorID,CanClone).
$ws.Range("K1").Value = "ActorID"
$ws.Range("J1").Value = "CanClone"

# Data rows.
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# New column J width (maps to xml width="14", no bestFit).
$ws.Columns("J").ColumnWidth = 13.285714285714286

# Page setup (paperSize=9 / A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to K9, matching the saved view state.
$ws.Range("K9").Select() | Out-Null
